$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "User Video" parameter row to become the "Mouse" parameter row.
$ws.Range("H1").Value = "Path to Mouse"
$ws.Range("H2").Value = "C:\Phase1\Backend\mouse_clicks"

# Update the active selection on the sheet to F2.
$ws.Range("F2").Select()
